$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "65.855.93"
Set-TextValue "E2" "  +0.17%  "

# Row 3
Set-TextValue "D3" "2.665.49"
Set-TextValue "E3" "  -0.49%  "

# Row 4
Set-TextValue "E4" "  +0.00%  "

# Row 5
Set-TextValue "D5" "598.74"
Set-TextValue "E5" "  -0.27%  "

# Row 6
Set-TextValue "D6" "158.66"
Set-TextValue "E6" "  +1.10%  "

# Row 7
Set-TextValue "E7" "  +4.95%  "

# Row 8
Set-TextValue "E8" "  +0.02%  "

# Row 9
Set-TextValue "E9" "  -2.94%  "

# Row 10
Set-TextValue "E10" "  +0.57%  "

# Row 11
Set-TextValue "E11" "  -0.50%  "

# Row 12
Set-TextValue "E12" "  +1.66%  "

# Row 13
Set-TextValue "D13" "29.08"
Set-TextValue "E13" "  -1.15%  "

# Row 14
Set-TextValue "E14" "  -2.17%  "

# Row 15
Set-TextValue "D15" "3.146.26"
Set-TextValue "E15" "  -0.49%  "

# Row 16
Set-TextValue "D16" "65.727.19"
Set-TextValue "E16" "  +0.15%  "

# Row 17
Set-TextValue "D17" "2.678.95"
Set-TextValue "E17" "  +0.07%  "

# Row 18
Set-TextValue "D18" "12.61"
Set-TextValue "E18" "  -2.30%  "

# Row 19
Set-TextValue "E19" "  +0.05%  "

# Row 20
Set-TextValue "D20" "7.51"
Set-TextValue "E20" "  -1.08%  "

# Row 21
Set-TextValue "D21" "351.77"
Set-TextValue "E21" "  -0.29%  "

# Row 22
Set-TextValue "E22" "  -0.07%  "

# Row 23
Set-TextValue "D23" "69.89"
Set-TextValue "E23" "  +0.11%  "

# Row 24
Set-TextValue "D24" "1.81"
Set-TextValue "E24" "  +10.61%  "

# Row 25
Set-TextValue "E25" "  -0.10%  "

# Row 26
Set-TextValue "D26" "9.66"
Set-TextValue "E26" "  -0.24%  "

# Row 27
Set-TextValue "E27" "  +1.49%  "

# Row 28
Set-TextValue "D28" "572.07"
Set-TextValue "E28" "  +7.92%  "

# Row 29
Set-TextValue "D29" "8.18"

# Row 30
Set-TextValue "E30" "  -2.55%  "

# Row 31
Set-TextValue "E31" "  -0.20%  "

# Row 32
Set-TextValue "E32" "  +0.85%  "

# Row 33
Set-TextValue "E33" "  +3.07%  "

# Row 34
Set-TextValue "D34" "6.75"
Set-TextValue "E34" "  +4.12%  "

# Row 35
Set-TextValue "D35" "5.57"
Set-TextValue "E35" "  +1.28%  "

# Row 36
Set-TextValue "E36" "  -0.13%  "

# Row 37
Set-TextValue "D37" "20.63"
Set-TextValue "E37" "  +0.09%  "

# Row 38
Set-TextValue "E38" "  -0.01%  "

# Row 39
Set-TextValue "E39" "  +0.65%  "

# Row 40
Set-TextValue "D40" "154.48"
Set-TextValue "E40" "  -2.32%  "

# Row 41
Set-TextValue "D41" "161.76"
Set-TextValue "E41" "  -1.90%  "

# Row 42
Set-TextValue "D42" "4.11"
Set-TextValue "E42" "  -0.80%  "

# Row 43
Set-TextValue "D43" "0.0620"
Set-TextValue "E43" "  +1.39%  "

# Row 44
Set-TextValue "E44" "  -0.64%  "

# Row 45
Set-TextValue "D45" "23.18"
Set-TextValue "E45" "  +1.15%  "

# Row 46
Set-TextValue "E46" "  +0.34%  "

# Row 47
Set-TextValue "D47" "0.0258"
Set-TextValue "E47" "  -0.50%  "

# Row 48
Set-TextValue "E48" "  +2.03%  "

# Row 49
Set-TextValue "D49" "19.82"
Set-TextValue "E49" "  -1.81%  "

# Row 50
Set-TextValue "E50" "  -6.79%  "

# Row 51
Set-TextValue "D51" "0.815"
Set-TextValue "E51" "  -0.47%  "
